{"js": "// \"Ch\u1ec9nh l\u1ea1i m\u1eabu 26\" \u2014 remove the leftover \"vnpt.SiteAddress\" placeholder\n// token that trails the \"\u0110\u1ecba ch\u1ec9: \" label in the B\u00ean A block, leaving the\n// label run untouched.\nconst body = context.document.body;\n\nconst results = body.search(\"vnpt.SiteAddress\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# \"Ch\u1ec9nh l\u1ea1i m\u1eabu 26\" \u2014 remove the leftover \"vnpt.SiteAddress\" placeholder\n# token that trails the \"\u0110\u1ecba ch\u1ec9: \" label in the B\u00ean A block, leaving the\n# label run untouched.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"vnpt.SiteAddress\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 0\n\nif ($find.Execute()) {\n    $range.Delete()\n}\n"}
